$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CALCULADORA")
$ws.Range("B3").Value = "772.37"
